$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 321, shifting the existing
# rows 321-349 down to 322-350 (dimension grows from A1:R349 to A1:R350).
$ws.Rows.Item(321).Insert()

# Populate the newly inserted row 321 with the new weekly record.
$ws.Range("A321").Value = 9
$ws.Range("B321").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C321").Value = "Metropolitana"
$ws.Range("D321").Value = 45106
$ws.Range("E321").Value = 13
$ws.Range("F321").Value = 100112026
$ws.Range("G321").Value = "Haba"
$ws.Range("H321").Value = "Sin especificar"
$ws.Range("I321").Value = "Primera"
$ws.Range("J321").Value = 70
$ws.Range("K321").Value = 16000
$ws.Range("L321").Value = 18000
$ws.Range("M321").Value = 17000
$ws.Range("N321").Value = "$/saco 25 kilos"
$ws.Range("O321").Value = "Provincia de Limarí"
$ws.Range("P321").Value = 680
$ws.Range("Q321").Value = 25
$ws.Range("R321").Value = "Hortaliza"
